# Update structure of excel:
#  - add a new VERSION column (Q) with header + "v2" for every data row
#  - add ENVIRONMENT ("klif") values for every data row (column O was blank)
#  - refresh RESPONSE TIME (column G) numbers
#  - change ERROR CODE (column K) from a bare newline to the literal "null"
#    and drop the wrap-text styling that column used to carry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new VERSION column -----------------------------------------------
# Give the new header cell (Q1) the same look as the other header cells
# (O1) by copying its formatting, then set the text.
$ws.Range("O1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Q1").Value = "VERSION"

# Match column Q's width to the other data columns (O is 35 chars wide).
$ws.Columns.Item(17).ColumnWidth = $ws.Columns.Item(15).ColumnWidth

# --- per-row data -------------------------------------------------------
$responseTimes = @("967","234","222","230","238","231","213","221","208","219","233","213","201","221","226","223","216","226","213","234","218")

for ($i = 0; $i -lt 21; $i++) {
    $row = $i + 2

    # RESPONSE TIME (column G) - keep it stored as text, same as before.
    $g = $ws.Cells.Item($row, 7)
    $g.NumberFormat = "@"
    $g.Value = $responseTimes[$i]
    $g.ClearFormats()

    # ERROR CODE (column K) - new literal text value, no more wrap-text style.
    $k = $ws.Cells.Item($row, 11)
    $k.Value = "null"
    $k.Style = "Normal"

    # ENVIRONMENT (column O) - now populated.
    $ws.Cells.Item($row, 15).Value = "klif"

    # VERSION (new column Q).
    $ws.Cells.Item($row, 17).Value = "v2"
}
